# Workbook and active worksheet references
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 - add Notes ("done")
$ws.Range("J10").Value = "done"

# Row 11 - add Names (Build/Testing = "hiếu") and Notes ("done")
$ws.Range("E11").Value = "hiếu"
$ws.Range("I11").Value = "hiếu"
$ws.Range("J11").Value = "done"

# Row 14 - add Notes ("done")
$ws.Range("J14").Value = "done"

# Row 15 - update Testing name to "quỳnh(hiếu)"
$ws.Range("I15").Value = "quỳnh(hiếu)"

# Row 16 - update Build/Testing names, add Notes
$ws.Range("E16").Value = "hoàng(hiếu)"
$ws.Range("I16").Value = "hoang (hiếu)"
$ws.Range("J16").Value = "done"

# Row 18 - update Build/Testing names, add Notes
$ws.Range("E18").Value = "quỳnh(hiếu)"
$ws.Range("I18").Value = "quỳnh(hiếu)"
$ws.Range("J18").Value = "done"

# Update the active selection to I23 (reflects author's last cursor position before saving)
$ws.Range("I23").Select()
